# Removes the trailing "Ver no Jupiter..." / copyright footer block that used
# to follow the "Requisitos" section (LOQ4046 requirement line), along with
# the blank separator paragraph right before it. The blank paragraph and the
# page-break paragraph that come after the footer block are left untouched.

$d = $word.ActiveDocument

$anchorText = "LOQ4046: Gestão da Produção e Logística (Requisito fraco)"

# Locate the anchor paragraph (the requirement line) by scanning the
# paragraphs collection, so the script is not tied to a fixed index.
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    # The three paragraphs right after the anchor are:
    #   anchorIndex + 1 -> blank separator paragraph
    #   anchorIndex + 2 -> "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   anchorIndex + 3 -> "© 2020 . Contact: ... Creative Commons Attribution"
    # All three (including their paragraph marks) get removed, while the
    # following blank paragraph + page-break paragraph are preserved.
    $startPara = $d.Paragraphs.Item($anchorIndex + 1)
    $endPara = $d.Paragraphs.Item($anchorIndex + 3)

    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
